$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cryptocurrency price/volume columns (D, E) store plain text values even
# when they look numeric (e.g. "0.690", "83.60"), since Excel would otherwise
# normalise them (dropping significant trailing zeros). For the handful of
# updated values that are unambiguous decimal numbers, force the cell to Text
# format first so the trailing zero / precision is preserved exactly as in the
# source data; every other cell is assigned directly.

$ws.Range('D2').Value = '66.925.65'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '3.103.66'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '576.45'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '177.52'
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.099.86'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '6.34'
$ws.Range('E10').Value = '  -2.39%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '36.13'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '3.621.66'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '66.902.75'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '3.103.08'
$ws.Range('D20').Value = '16.73'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').Value = '479.56'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').Value = '7.78'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.690'
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.60'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = '12.57'
$ws.Range('E25').Value = '  -3.63%  '
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('D27').Value = '10.11'
$ws.Range('E27').Value = '  -3.79%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '7.91'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.30'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('D31').Value = '2.61'
$ws.Range('E31').Value = '  -2.05%  '
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').Value = '0.0₃0941'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '48.33'
$ws.Range('E36').Value = '  +3.23%  '
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('D38').Value = '0.941'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').Value = '0.312'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').Value = '8.33'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').Value = '2.71'
$ws.Range('E44').Value = '  +5.85%  '
$ws.Range('D45').Value = '2.797.22'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '372.53'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.80'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0344'
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D50').Value = '25.76'
$ws.Range('E50').Value = '  +3.37%  '
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').Value = '  +2.49%  '
